$wb = $excel.ActiveWorkbook
$wsProgress = $wb.Worksheets.Item("Progress")
$wsJournal  = $wb.Worksheets.Item("Journal")

# --- Progress sheet: update row 29 (Link PC(Simulink) to Board(C program)) ---
# H29 = Start Date, I29 = End Date, J29 = Status
$wsProgress.Range("H29").Value2 = 45806
$wsProgress.Range("I29").Value2 = 45814
$wsProgress.Range("J29").Value2 = "done"

# --- Journal sheet: append new journal entry row 24 ---
$wsJournal.Range("A23:C23").Copy($wsJournal.Range("A24:C24"))
$wsJournal.Range("A24").Value2 = "Link PC(Simulink) to Board(C program)"
$wsJournal.Range("B24").Value2 = 45814
$journalText = "Linking the board went smoothly; I connected it via Ethernet using the eth0 interface, brought it up with ifconfig eth0 up, and obtained an IP address through DHCP using udhcpc -i eth0. I then connected to the assigned IP from the MATLAB dataset interface for data transmission.`nWhile testing Realtime Dataset Mode on the DE1-SoC board using the setup described above, a segmentation fault occurred. The root cause was likely a stack overflow due to oversized buffers declared within the processing_pipeline function.`nOn many embedded Linux systems, including DE1-SoC, the default pthread stack size is often only 64 KB or 128 KB(I tried confirming this from the board terminal, I couldnt). Although I could have increased it using pthread_attr_setstacksize(), I decided against it since future implementations will run on memory-constrained environments like bare-metal systems. `nInstead, I refactored all large buffer variables inside the pipeline to be file-scoped static variables, allowing them to reside in the .bss segment (static memory) instead of the thread’s stack. This fixed the issue. "
$wsJournal.Range("C24").Value2 = $journalText
$wsJournal.Range("A24:C24").EntireRow.AutoFit() | Out-Null
$wsJournal.Rows.Item(24).RowHeight = 145.75

# --- View state updates ---
# Progress sheet keeps its scroll position but loses tabSelected & the active cell moves to K30
$wsProgress.Activate()
$wsProgress.Range("K30").Select() | Out-Null
$winProgress = $excel.ActiveWindow
$winProgress.ScrollRow = 24
$winProgress.ScrollColumn = 1

# Journal sheet becomes the active/selected tab, with selection on the newly added C24 cell
$wsJournal.Activate()
$wsJournal.Range("C24").Select() | Out-Null
$winJournal = $excel.ActiveWindow
$winJournal.ScrollRow = 23
$winJournal.ScrollColumn = 1
